$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 804, shifting rows 804:845 down to 805:846
$ws.Rows.Item(804).Insert()

# Populate the newly inserted row 804
# Force column A to remain plain text (matching the rest of the sheet)
# rather than being auto-converted into a date serial number.
$ws.Cells.Item(804, 1).NumberFormat = "@"
$ws.Cells.Item(804, 1).Value = "2026/02/16"
$ws.Cells.Item(804, 1).Style = "Normal"
$ws.Cells.Item(804, 2).Value = "月"
$ws.Cells.Item(804, 3).Value = 13
$ws.Cells.Item(804, 4).Value = 65
